$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns at R:T (shifts old R..AE -> U..AH),
# carrying the bold/bordered header style into the new blank header cells.
$ws.Range("R1:T1").EntireColumn.Insert()

# --- Row 1 (header) : fill the 3 newly inserted header cells ---
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# --- Row 2 (data) : fill previously-blank "unknown" placeholder cells ---
$ws.Range("D2").Value = "unknown"
$ws.Range("E2").Value = "unknown"
$ws.Range("F2").Value = "unknown"
$ws.Range("G2").Value = "unknown"
$ws.Range("H2").Value = "unknown"
$ws.Range("I2").Value = "unknown"
$ws.Range("J2").Value = "unknown"

# Values for the 3 newly-inserted data cells (shifted columns keep their
# existing values automatically via the column insert above).
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 3
$ws.Range("T2").Value = 0

# The old "arts" column (now U2) changes representation from boolean to
# a plain number.
$ws.Range("U2").Value = 0

Write-Host "done"
